# Phase1-answer/Question1/Q2Data.xlsx — "debugs input variebles of first part. all the answers changed"
#
# Adds two new input rows (vegrefinelimit=200, oilrefinelimit=250) to the
# "input" sheet, and updates the "output" sheet with the recomputed solver
# results that follow from that change (literal values pasted in by the
# author, plus the SUMPRODUCT-based weighted-average formulas in M12:R12
# which recalc automatically from the new E12:J16 figures).

$wb = $excel.ActiveWorkbook
$wsInput  = $wb.Worksheets.Item("input")
$wsOutput = $wb.Worksheets.Item("output")

# --- input sheet: two new rows, formatted like the existing B9:C12 block ---
$wsInput.Range("B12").Copy()
$wsInput.Range("B13").PasteSpecial(-4122) # xlPasteFormats
$wsInput.Range("B13").Value = "vegrefinelimit"
$wsInput.Range("C13").Value = 200

$wsInput.Range("B12").Copy()
$wsInput.Range("B14").PasteSpecial(-4122) # xlPasteFormats
$wsInput.Range("B14").Value = "oilrefinelimit"
$wsInput.Range("C14").Value = 250

# --- output sheet: updated solver results ---
$wsOutput.Range("M3").Value = 405000
$wsOutput.Range("I4").Value = 199.99999999999994
$wsOutput.Range("J4").Value = 532.63888888888891
$wsOutput.Range("M4").Value = 246493.14236111112
$wsOutput.Range("J5").Value = 467.36111111111114
$wsOutput.Range("M5").Value = 54119.531250000007
$wsOutput.Range("F6").Value = 32.5
$wsOutput.Range("M6").Value = 104387.32638888888
$wsOutput.Range("F7").Value = 37.500000000000341
$wsOutput.Range("J7").Value = 732.03125
$wsOutput.Range("G8").Value = 697.96874999999977

$wsOutput.Range("E12").Value = 159.25925925925921
$wsOutput.Range("F12").Value = 127.77777777777777
$wsOutput.Range("G12").Value = 106.48148148148151
$wsOutput.Range("H12").Value = 53.240740740740755
$wsOutput.Range("I12").Value = 126.62037037037035
$wsOutput.Range("J12").Value = 159.25925925925921

$wsOutput.Range("E13").Value = 40.74074074074079
$wsOutput.Range("F13").Value = 72.222222222222229
$wsOutput.Range("G13").Value = 93.518518518518491
$wsOutput.Range("H13").Value = 146.75925925925924
$wsOutput.Range("I13").Value = 73.379629629629648
$wsOutput.Range("J13").Value = 40.74074074074079

$wsOutput.Range("I14").Value = 32.5

$wsOutput.Range("E15").Value = 250
$wsOutput.Range("F15").Value = 143.75000000000017
$wsOutput.Range("G15").Value = 71.875000000000085
$wsOutput.Range("H15").Value = 35.937500000000043
$wsOutput.Range("I15").Value = 17.968750000000021
$wsOutput.Range("J15").Value = 250

$wsOutput.Range("F16").Value = 106.24999999999984
$wsOutput.Range("G16").Value = 178.12499999999991
$wsOutput.Range("H16").Value = 214.06249999999994
$wsOutput.Range("I16").Value = 199.53124999999997

$wsOutput.Range("E20").Value = 340.74074074074076
$wsOutput.Range("F20").Value = 212.96296296296302
$wsOutput.Range("G20").Value = 106.48148148148151
$wsOutput.Range("H20").Value = 53.240740740740755
$wsOutput.Range("I20").Value = 126.62037037037035

$wsOutput.Range("E21").Value = 459.25925925925924
$wsOutput.Range("F21").Value = 387.03703703703701
$wsOutput.Range("G21").Value = 293.51851851851853
$wsOutput.Range("H21").Value = 146.7592592592593
$wsOutput.Range("I21").Value = 73.379629629629648

$wsOutput.Range("F22").Value = 532.5
$wsOutput.Range("G22").Value = 532.5
$wsOutput.Range("H22").Value = 532.5

$wsOutput.Range("E23").Value = 250
$wsOutput.Range("F23").Value = 143.75000000000017
$wsOutput.Range("G23").Value = 71.875000000000085
$wsOutput.Range("H23").Value = 35.937500000000043
$wsOutput.Range("I23").Value = 17.968750000000021

$wsOutput.Range("F24").Value = 393.75000000000017
$wsOutput.Range("G24").Value = 913.59375
$wsOutput.Range("H24").Value = 699.53125

# --- restore the selections the author left in each sheet ---
$wsInput.Activate()
$wsInput.Range("F10").Select()

$wsOutput.Activate()
$wsOutput.Range("H12").Select()
